# Auto-update draw results: append the 2025-10-17 Pick 4 row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1   # xlUp = -4162; first empty row below the data

$values = @{
    1 = "2025-10-17"
    2 = "Pick 4"
    3 = "251017"
    4 = "7-5-6-2"
    5 = "2025-10-17T21:37:16.657+04:00"
}

foreach ($col in 1..5) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text storage (these look like a date / a number to Excel's
    # auto-detection) while keeping the cell's number format as "General",
    # matching every other data row in the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col]
    $cell.Style = "Normal"
}

# Keep the "numbers stored as text" warning suppressed over the data range,
# extending it to cover the newly-added row (mirrors Excel's own behavior
# when you dismiss the green-triangle warning on the new cells).
try {
    $lastRow = $row
    $ws.Range("A1:E$lastRow").Errors.Item(9).Ignore = $true
} catch {
    # Older/limited COM surfaces may not expose per-error Ignore; safe to skip.
}
